$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.280.58"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.638.33"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "528.42"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "144.93"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "0.104"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "3.106.30"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "59.252.51"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "20.93"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.691.69"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "342.67"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "4.46"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "10.61"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "65.69"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").Value = "0.418"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "7.24"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "6.45"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").Value = "19.06"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").Value = "150.17"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "4.19"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "0.886"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "0.866"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "36.64"
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "3.66"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.604"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.0975"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "269.53"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "19.44"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").Value = "0.0538"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "2.041.56"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0230"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.69"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").Value = "18.81"
$ws.Range("E51").Value = "  -1.31%  "
